$d = $word.ActiveDocument

$body = '<w:p><w:r><w:t>Logo + Tagline</w:t></w:r><w:r><w:tab/></w:r></w:p>' +
        '<w:p><w:r><w:t>Login</w:t></w:r><w:r><w:t>/Register</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t xml:space="preserve"> Button</w:t></w:r></w:p>' +
        '<w:p><w:r><w:t>Main Image</w:t></w:r><w:r><w:t xml:space="preserve"> With </w:t></w:r>' +
        '<w:r><w:t>A-Head</w:t></w:r><w:r><w:t xml:space="preserve"> &amp; </w:t></w:r>' +
        '<w:r><w:t>B-Head</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
        '<w:p><w:r><w:t xml:space="preserve">Alternating </w:t></w:r><w:r><w:t>Sub Images</w:t></w:r></w:p>' +
        '<w:p><w:r><w:t>Menu Button -&gt; Animated Pop out Side Menu</w:t></w:r></w:p>' +
        '<w:p><w:r><w:t>Footer</w:t></w:r></w:p>' +
        '<w:p/><w:p/><w:p/>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
